$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Masthead: issue number + week-covered dates ---
$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# --- Weekly Crime Complaints table refresh (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = -11.111111111111
$ws.Range("I14").Value = 54
$ws.Range("J14").Value = 61
$ws.Range("K14").Value = -11.475409836065
$ws.Range("L14").Value = -28.947368421052
$ws.Range("M14").Value = -15.625
$ws.Range("N14").Value = -79.622641509434

# Row 15
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 200
$ws.Range("F15").Value = 45
$ws.Range("G15").Value = 32
$ws.Range("H15").Value = 40.625
$ws.Range("I15").Value = 284
$ws.Range("J15").Value = 228
$ws.Range("K15").Value = 24.561403508771
$ws.Range("L15").Value = 33.962264150943
$ws.Range("M15").Value = 89.333333333333
$ws.Range("N15").Value = -24.867724867724

# Row 16
$ws.Range("C16").Value = 104
$ws.Range("D16").Value = 117
$ws.Range("E16").Value = -11.111111111111
$ws.Range("F16").Value = 417
$ws.Range("G16").Value = 431
$ws.Range("H16").Value = -3.248259860788
$ws.Range("I16").Value = 2492
$ws.Range("J16").Value = 2602
$ws.Range("K16").Value = -4.227517294388
$ws.Range("L16").Value = 1.465798045602
$ws.Range("M16").Value = 11.449016100178
$ws.Range("N16").Value = -70.522829429855

# Row 17
$ws.Range("C17").Value = 149
$ws.Range("D17").Value = 202
$ws.Range("E17").Value = -26.237623762376
$ws.Range("F17").Value = 780
$ws.Range("G17").Value = 807
$ws.Range("H17").Value = -3.345724907063
$ws.Range("I17").Value = 4623
$ws.Range("J17").Value = 4415
$ws.Range("K17").Value = 4.711211778029
$ws.Range("L17").Value = 9.549763033175
$ws.Range("M17").Value = 99.870298313878
$ws.Range("N17").Value = -3.767693588676

# Row 18
$ws.Range("C18").Value = 51
$ws.Range("D18").Value = 57
$ws.Range("E18").Value = -10.526315789473
$ws.Range("F18").Value = 223
$ws.Range("G18").Value = 241
$ws.Range("H18").Value = -7.468879668049
$ws.Range("I18").Value = 1514
$ws.Range("J18").Value = 1544
$ws.Range("K18").Value = -1.943005181347
$ws.Range("L18").Value = -4.95919648462
$ws.Range("M18").Value = -10.36116044997
$ws.Range("N18").Value = -84.713247172859

# Row 19
$ws.Range("C19").Value = 227
$ws.Range("D19").Value = 218
$ws.Range("E19").Value = 4.128440366972
$ws.Range("F19").Value = 748
$ws.Range("G19").Value = 716
$ws.Range("H19").Value = 4.469273743016
$ws.Range("I19").Value = 4805
$ws.Range("J19").Value = 4803
$ws.Range("K19").Value = 0.041640641265
$ws.Range("L19").Value = 16.428398352314
$ws.Range("M19").Value = 104.294217687075
$ws.Range("N19").Value = 25.228042741725

# Row 20
$ws.Range("C20").Value = 96
$ws.Range("D20").Value = 83
$ws.Range("E20").Value = 15.662650602409
$ws.Range("F20").Value = 367
$ws.Range("G20").Value = 345
$ws.Range("H20").Value = 6.376811594202
$ws.Range("I20").Value = 2377
$ws.Range("J20").Value = 2163
$ws.Range("K20").Value = 9.893666204345
$ws.Range("L20").Value = -15.768958185683
$ws.Range("M20").Value = 122.983114446529
$ws.Range("N20").Value = -70.578041836861

# Row 21
$ws.Range("C21").Value = 641
$ws.Range("D21").Value = 682
$ws.Range("E21").Value = -6.011730205278
$ws.Range("F21").Value = 2588
$ws.Range("G21").Value = 2581
$ws.Range("H21").Value = 0.271212708252
$ws.Range("I21").Value = 16149
$ws.Range("J21").Value = 15816
$ws.Range("K21").Value = 2.105462822458
$ws.Range("L21").Value = 4.14678189088
$ws.Range("M21").Value = 63.617021276595
$ws.Range("N21").Value = -54.791299235743

# Row 22
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 40
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 20
$ws.Range("H22").Value = -10
$ws.Range("I22").Value = 161
$ws.Range("J22").Value = 177
$ws.Range("K22").Value = -9.039548022598
$ws.Range("L22").Value = 1.898734177215
$ws.Range("M22").Value = -5.847953216374

# Row 23
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = -31.25
$ws.Range("F23").Value = 132
$ws.Range("G23").Value = 142
$ws.Range("H23").Value = -7.042253521126
$ws.Range("I23").Value = 828
$ws.Range("J23").Value = 906
$ws.Range("K23").Value = -8.609271523178
$ws.Range("L23").Value = -12.658227848101
$ws.Range("M23").Value = 50.2722323049

# Row 24
$ws.Range("C24").Value = 398
$ws.Range("D24").Value = 344
$ws.Range("E24").Value = 15.697674418604
$ws.Range("F24").Value = 1495
$ws.Range("G24").Value = 1141
$ws.Range("H24").Value = 31.025416301489
$ws.Range("I24").Value = 9501
$ws.Range("J24").Value = 8474
$ws.Range("K24").Value = 12.11942412084
$ws.Range("L24").Value = 0.560965283657
$ws.Range("M24").Value = 44.942791762013

# Row 25
$ws.Range("C25").Value = 145
$ws.Range("D25").Value = 148
$ws.Range("E25").Value = -2.027027027027
$ws.Range("F25").Value = 529
$ws.Range("G25").Value = 452
$ws.Range("H25").Value = 17.035398230088
$ws.Range("I25").Value = 3157
$ws.Range("J25").Value = 3388
$ws.Range("K25").Value = -6.818181818181
$ws.Range("L25").Value = -23.112518265952

# Row 26
$ws.Range("C26").Value = 253
$ws.Range("D26").Value = 242
$ws.Range("E26").Value = 4.545454545454
$ws.Range("F26").Value = 995
$ws.Range("G26").Value = 985
$ws.Range("H26").Value = 1.015228426395
$ws.Range("I26").Value = 5936
$ws.Range("J26").Value = 5924
$ws.Range("K26").Value = 0.202565833896
$ws.Range("L26").Value = 6.379928315412
$ws.Range("M26").Value = 0.202565833896

# Row 27
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 54
$ws.Range("G27").Value = 52
$ws.Range("H27").Value = 3.846153846153
$ws.Range("I27").Value = 357
$ws.Range("J27").Value = 355
$ws.Range("K27").Value = 0.56338028169
$ws.Range("L27").Value = 2.292263610315

# Row 28
$ws.Range("C28").Value = 26
$ws.Range("D28").Value = 18
$ws.Range("E28").Value = 44.444444444444
$ws.Range("F28").Value = 83
$ws.Range("G28").Value = 86
$ws.Range("H28").Value = -3.488372093023
$ws.Range("I28").Value = 587
$ws.Range("J28").Value = 639
$ws.Range("K28").Value = -8.137715179968
$ws.Range("L28").Value = 5.197132616487

# Row 29
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = 9
$ws.Range("E29").Value = -11.111111111111
$ws.Range("F29").Value = 30
$ws.Range("G29").Value = 46
$ws.Range("H29").Value = -34.782608695652
$ws.Range("I29").Value = 161
$ws.Range("J29").Value = 207
$ws.Range("K29").Value = -22.222222222222
$ws.Range("L29").Value = -25.462962962963
$ws.Range("M29").Value = -31.779661016949
$ws.Range("N29").Value = -77.065527065527

# Row 30
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = 16.666666666666
$ws.Range("F30").Value = 28
$ws.Range("G30").Value = 37
$ws.Range("H30").Value = -24.324324324324
$ws.Range("I30").Value = 141
$ws.Range("J30").Value = 168
$ws.Range("K30").Value = -16.071428571428
$ws.Range("L30").Value = -19.886363636363
$ws.Range("M30").Value = -28.426395939086
$ws.Range("N30").Value = -77.725118483412

# Row 14, col C: was text "0", now a real number
$ws.Range("C14").Value = 2

# Row 33: mix of numeric <-> text ("0" / "***.*") changes plus a couple of numbers
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "***.*"
$ws.Range("I33").Value = 14
$ws.Range("K33").Value = -41.666666666666
$ws.Range("L33").Value = -41.666666666666

# --- Column E narrows back to the shared 6.168446-char best-fit width used by columns C/D/F/G/H ---
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth
